# Regulation17 template update:
#  - Rename Sheet1 -> "Species Information"
#  - Add a new sheet "More Species Information" after it
#  - Bold the header rows on both sheets
#  - Add sample data rows (Frog/Kangaroo on sheet 1, Rat/Dog on sheet 2)
#  - Widen column B on sheet 1 for the longer "Species Count" header
#  - Make "More Species Information" the active/selected sheet

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Species Information" ------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Species Information"

# Header row (bold)
$ws1.Range("A1:D1").Font.Bold = $true

# Data rows
$ws1.Range("A2").Value = "Frog"
$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = 35
$ws1.Range("D2").Value = 170

$ws1.Range("A3").Value = "Kangaroo"
$ws1.Range("B3").Value = 2
$ws1.Range("C3").Value = 35
$ws1.Range("D3").Value = 170

# Column widths: A stays near-default, B widened for "Species Count"
$ws1.Columns.Item(1).ColumnWidth = 10.5
$ws1.Columns.Item(2).ColumnWidth = 28.5

# --- Sheet 2: "More Species Information" --------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "More Species Information"

$ws2.Range("A1").Value = "Species"
$ws2.Range("B1").Value = "Species Count"
$ws2.Range("A1:B1").Font.Bold = $true

$ws2.Range("A2").Value = "Rat"
$ws2.Range("B2").Value = 3

$ws2.Range("A3").Value = "Dog"
$ws2.Range("B3").Value = 4

# Match the printed page setup used on sheet 1 (margins, paper size,
# orientation and header/footer text)
$ps2 = $ws2.PageSetup
$ps2.LeftMargin = 56.699999999999996
$ps2.RightMargin = 56.699999999999996
$ps2.TopMargin = 75.80000000000017
$ps2.BottomMargin = 75.80000000000017
$ps2.HeaderMargin = 56.699999999999996
$ps2.FooterMargin = 56.699999999999996
$ps2.PaperSize = 9
$ps2.Orientation = 1
$ps2.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ps2.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

# --- Selections / active views ------------------------------------------
$ws1.Range("C12").Select() | Out-Null
$ws2.Range("D16").Select() | Out-Null

# Make the new sheet the active tab
$ws2.Activate() | Out-Null
